$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.015.86"
$ws.Range("E2").Value = "  +2.33%  "

# Row 3
$ws.Range("D3").Value = "3.187.39"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.37"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.91"
$ws.Range("E6").Value = "  +3.85%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -2.14%  "

# Row 9
$ws.Range("E9").Value = "  +0.49%  "

# Row 10
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.431"
$ws.Range("E11").Value = "  -1.85%  "

# Row 12
$ws.Range("D12").Value = "3.736.73"
$ws.Range("E12").Value = "  +1.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.138"
$ws.Range("E13").Value = "  -2.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.77"
$ws.Range("E14").Value = "  -0.60%  "

# Row 15
$ws.Range("E15").Value = "  +0.02%  "

# Row 16
$ws.Range("D16").Value = "60.035.28"
$ws.Range("E16").Value = "  +2.28%  "

# Row 17
$ws.Range("D17").Value = "3.200.11"
$ws.Range("E17").Value = "  +1.46%  "

# Row 18
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +1.82%  "

# Row 20
$ws.Range("E20").Value = "  +0.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.78"
$ws.Range("E21").Value = "  -1.87%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("E23").Value = "  -1.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.44"
$ws.Range("E24").Value = "  -0.22%  "

# Row 25
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.59"
$ws.Range("E26").Value = "  +4.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  -0.84%  "

# Row 28
$ws.Range("E28").Value = "  +0.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.44"
$ws.Range("E29").Value = "  +0.39%  "

# Row 30
$ws.Range("E30").Value = "  +0.54%  "

# Row 31
$ws.Range("E31").Value = "  +0.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  +2.79%  "

# Row 33
$ws.Range("E33").Value = "  +5.02%  "

# Row 34
$ws.Range("E34").Value = "  +3.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.95"
$ws.Range("E35").Value = "  -0.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.35"
$ws.Range("E36").Value = "  +1.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.39"
$ws.Range("E37").Value = "  +5.77%  "

# Row 38
$ws.Range("D38").Value = "2.782.40"
$ws.Range("E38").Value = "  +5.16%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0708"
$ws.Range("E39").Value = "  +3.25%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  +0.47%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0307"
$ws.Range("E41").Value = "  +6.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.21"
$ws.Range("E42").Value = "  -1.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.82"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("E45").Value = "  +1.09%  "

# Row 46
$ws.Range("D46").Value = "3.228.64"
$ws.Range("E46").Value = "  +1.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.982"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("E49").Value = "  +5.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.56"
$ws.Range("E50").Value = "  +2.83%  "

# Row 51
$ws.Range("E51").Value = "  +0.05%  "
